$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column map (1-based) ---
# 1=A Date  2=B Cumul cas positifs(formula)  3=C Nb nouveaux cas positifs
# 4=D Nb nouvelles admissions a l'hopital    5=E Patients SI total
# 6=F Patients intubes                       7=G Patients hosp hors SI
# 8=H Total hospitalisations(formula)        9=I Nb nouvelles sorties
# 10=J Cumul deces(formula)                  11=K Nb nouveaux deces(formula)
# 12=L Nb nouveaux deces hopital             13=M Nb nouveaux deces extra-hosp
# 14=N cas isolement  15=O contacts quarantaine  16=P voyageurs quarantaine

# --- Simple numeric corrections on existing rows (columns C,D,G,N,O,P only;
#     B/H/J/K are formulas and recompute automatically) ---

$ws.Cells.Item(199, 15).Value = 225

$ws.Cells.Item(200, 15).Value = 231

$ws.Cells.Item(201, 15).Value = 254

$ws.Cells.Item(202, 15).Value = 293

$ws.Cells.Item(205, 15).Value = 383

$ws.Cells.Item(208, 15).Value = 417
$ws.Cells.Item(208, 16).Value = 440

$ws.Cells.Item(209, 14).Value = 113
$ws.Cells.Item(209, 15).Value = 382
$ws.Cells.Item(209, 16).Value = 462

$ws.Cells.Item(210, 14).Value = 118
$ws.Cells.Item(210, 15).Value = 396
$ws.Cells.Item(210, 16).Value = 449

$ws.Cells.Item(211, 14).Value = 116
$ws.Cells.Item(211, 15).Value = 400
$ws.Cells.Item(211, 16).Value = 439

$ws.Cells.Item(212, 3).Value = 10
$ws.Cells.Item(212, 14).Value = 109
$ws.Cells.Item(212, 15).Value = 423
$ws.Cells.Item(212, 16).Value = 462

$ws.Cells.Item(213, 3).Value = 6
$ws.Cells.Item(213, 4).Value = 1
$ws.Cells.Item(213, 7).Value = 9
$ws.Cells.Item(213, 14).Value = 107
$ws.Cells.Item(213, 15).Value = 442
$ws.Cells.Item(213, 16).Value = 434

# --- Newly-entered rows 214-216 (previously blank input columns) ---

$ws.Cells.Item(214, 3).Value = 5
$ws.Cells.Item(214, 4).Value = 0
$ws.Cells.Item(214, 5).Value = 1
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 8
$ws.Cells.Item(214, 9).Value = 0
$ws.Cells.Item(214, 12).Value = "0"
$ws.Cells.Item(214, 13).Value = "0"
$ws.Cells.Item(214, 14).Value = 94
$ws.Cells.Item(214, 15).Value = 389
$ws.Cells.Item(214, 16).Value = 393

$ws.Cells.Item(215, 3).Value = 2
$ws.Cells.Item(215, 4).Value = 0
$ws.Cells.Item(215, 5).Value = 1
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 8
$ws.Cells.Item(215, 9).Value = 0
$ws.Cells.Item(215, 12).Value = "0"
$ws.Cells.Item(215, 13).Value = "0"
$ws.Cells.Item(215, 14).Value = 91
$ws.Cells.Item(215, 15).Value = 346
$ws.Cells.Item(215, 16).Value = 359

$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 0
$ws.Cells.Item(216, 5).Value = 1
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 8
$ws.Cells.Item(216, 9).Value = 0
$ws.Cells.Item(216, 12).Value = "0"
$ws.Cells.Item(216, 13).Value = "0"
$ws.Cells.Item(216, 14).Value = 74
$ws.Cells.Item(216, 15).Value = 299
$ws.Cells.Item(216, 16).Value = 304

# --- View state: update the selected/active cell in the frozen bottom-right
#     pane (matches the author re-opening the file further down the sheet) ---
$ws.Range("P125").Select()
